$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values per row (matching the TPM update described in the diff)
$data = @{
    2 = @{ G = 0.2284785; H = 0.456957; M = 5.978421000000001; N = 17.935263;
           O = 0.05704457007880161; P = 0.06242884486533885; Q = 1.3659406624485;
           R = 8.195643974691; S = 0.05704457007880161; T = 0.06242884486533885 }
    3 = @{ G = 0.2284785; H = 0.456957;
           O = 0.6646576013185088; P = 0.7273927426214574; Q = 15.9153245084025;
           R = 95.49194705041499; S = 0.6646576013185088; T = 0.7273927426214574 }
    4 = @{ G = 0.2284785; H = 0.456957; M = 1.290243; N = 3.870729;
           O = 0.01231116999491725; P = 0.01347318632889677; Q = 0.2947927852755;
           R = 1.768756711653; S = 0.01231116999491725; T = 0.01347318632889677 }
    5 = @{ G = 0.2284785; H = 0.456957; M = 27.1166075; N = 54.233215;
           O = 0.2587397603536297; P = 0.1887743138075849; Q = 6.19556180668875;
           R = 24.782247226755; S = 0.2587397603536297; T = 0.1887743138075849 }
    6 = @{ G = 0.2284785; H = 0.456957; M = 0.759494; N = 2.278482;
           O = 0.00724689825414258; P = 0.007930912376722157; Q = 0.173528049879;
           R = 1.041168299274; S = 0.00724689825414258; T = 0.007930912376722157 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
